$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Coles / SJ300 (rename "coles" -> "Coles_SJ300", condition unchanged)
$ws.Range("A2").Value = "Coles_SJ300"
$ws.Range("B2").Value = "specialpricecat product"
$ws.Range("C2").Value = "specialpricecat==012 & product==SJ300"

# Row 3: Woolworths / SJ300
$ws.Range("A3").Value = "Woolworths_SJ300"
$ws.Range("B3").Value = "specialpricecat product"
$ws.Range("C3").Value = "specialpricecat==010 & product==SJ300 "

# Row 4: Indies / SJ300 (replaces old "shop" row)
$ws.Range("A4").Value = "Indies_SJ300"
$ws.Range("B4").Value = "product"
$ws.Range("C4").Value = "(specialpricecat==088 | specialpricecat==028 | specialpricecat==038 | specialpricecat==048 | specialpricecat==058 | specialpricecat==068 | specialpricecat==078 ) & product==SJ300"

# Row 5: Distributors / SJ300 (replaces old "indies" row)
$ws.Range("A5").Value = "Distributors_SJ300"
$ws.Range("B5").Value = "code"
$ws.Range("C5").Value = "(specialpricecat==080 | specialpricecat==020 | specialpricecat==030 | specialpricecat==040 | specialpricecat==050 | specialpricecat==060 | specialpricecat==070 ) & product==SJ300"

# Row 6: Coles / RJ300 (replaces old "distributors" row)
$ws.Range("A6").Value = "Coles_RJ300"
$ws.Range("B6").Value = "specialpricecat product"
$ws.Range("C6").Value = "specialpricecat==012 & product==RJ300"

# Row 7: Woolworths / RJ300 (new row)
$ws.Range("A7").Value = "Woolworths_RJ300"
$ws.Range("B7").Value = "specialpricecat product"
$ws.Range("C7").Value = "specialpricecat==010 & product==RJ300 "

# Row 8: Indies / RJ300 (new row)
$ws.Range("A8").Value = "Indies_RJ300"
$ws.Range("B8").Value = "product"
$ws.Range("C8").Value = "(specialpricecat==088 | specialpricecat==028 | specialpricecat==038 | specialpricecat==048 | specialpricecat==058 | specialpricecat==068 | specialpricecat==078 ) & product==RJ300"

# Row 9: Distributors / RJ300 (new row)
$ws.Range("A9").Value = "Distributors_RJ300"
$ws.Range("B9").Value = "code"
$ws.Range("C9").Value = "(specialpricecat==080 | specialpricecat==020 | specialpricecat==030 | specialpricecat==040 | specialpricecat==050 | specialpricecat==060 | specialpricecat==070 ) & product==RJ300"

# Row 10: Indies / CFJ300 (new row)
$ws.Range("A10").Value = "Indies_CFJ300"
$ws.Range("B10").Value = "product"
$ws.Range("C10").Value = "(specialpricecat==088 | specialpricecat==028 | specialpricecat==038 | specialpricecat==048 | specialpricecat==058 | specialpricecat==068 | specialpricecat==078 ) & product==CFJ300"

# Row 11: Distributors / CFJ300 (new row)
$ws.Range("A11").Value = "Distributors_CFJ300"
$ws.Range("B11").Value = "code"
$ws.Range("C11").Value = "(specialpricecat==080 | specialpricecat==020 | specialpricecat==030 | specialpricecat==040 | specialpricecat==050 | specialpricecat==060 | specialpricecat==070 ) & product==CFJ300"

# Widen column A to fit the longer labels (24.45 chars; COM storage quantizes to
# the nearest pixel, so 23.6 is the input that lands closest to 24.45 on save)
$ws.Columns.Item(1).ColumnWidth = 23.6

# Match the saved selection state (B10:B11 selected, active cell B10)
$ws.Range("B10:B11").Select()
